$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 1-9 are untouched by this edit. Rebuild rows 10-25 from scratch so that
# leftover styles/row-heights from the old layout do not linger.
$ws.Range("A10:C25").EntireRow.Delete()

# Row 10
$ws.Range('A9').Copy()
$ws.Range('A10').PasteSpecial(-4122)
$ws.Range('A10').Value = 'Objetivos:'
$ws.Range('B9').Copy()
$ws.Range('B10').PasteSpecial(-4122)
$ws.Range('B10').Value = '3444370 - Rita de Cássia Lacerda Brambilla Rodrigues'
$ws.Range('C9').Copy()
$ws.Range('C10').PasteSpecial(-4122)
$ws.Range('C10').Value = '3444370 - Rita de Cássia Lacerda Brambilla Rodrigues'
$ws.Rows.Item(10).RowHeight = 60

# Row 11
$ws.Range('A9').Copy()
$ws.Range('A11').PasteSpecial(-4122)
$ws.Range('A11').Value = 'Objectives:'
$ws.Rows.Item(11).RowHeight = 60

# Row 12
$ws.Range('A9').Copy()
$ws.Range('A12').PasteSpecial(-4122)
$ws.Range('A12').Value = 'Docentes responsáveis:'

# Row 13
$ws.Range('A9').Copy()
$ws.Range('A13').PasteSpecial(-4122)
$ws.Range('A13').Value = 'Programa resumido:'
$ws.Range('B9').Copy()
$ws.Range('B13').PasteSpecial(-4122)
$ws.Range('B13').Value = 'Semestral'
$ws.Range('C9').Copy()
$ws.Range('C13').PasteSpecial(-4122)
$ws.Range('C13').Value = 'Semestral'
$ws.Rows.Item(13).RowHeight = 60

# Row 14
$ws.Range('A9').Copy()
$ws.Range('A14').PasteSpecial(-4122)
$ws.Range('A14').Value = 'Short syllabus:'
$ws.Range('B9').Copy()
$ws.Range('B14').PasteSpecial(-4122)
$ws.Range('B14').Value = 'Proceedings of biochemical processes since the steps of medium preparation and sterilization until the steps of products recuperation and characterization.'
$ws.Range('C9').Copy()
$ws.Range('C14').PasteSpecial(-4122)
$ws.Range('C14').Value = 'Proceedings of biochemical processes since the steps of medium preparation and sterilization until the steps of products recuperation and characterization.'
$ws.Rows.Item(14).RowHeight = 60

# Row 15
$ws.Range('A9').Copy()
$ws.Range('A15').PasteSpecial(-4122)
$ws.Range('A15').Value = 'Programa:'
$ws.Range('B9').Copy()
$ws.Range('B15').PasteSpecial(-4122)
$ws.Range('B15').Value = '01/01/2012'
$ws.Range('C9').Copy()
$ws.Range('C15').PasteSpecial(-4122)
$ws.Range('C15').Value = '01/01/2012'
$ws.Rows.Item(15).RowHeight = 120

# Row 16
$ws.Range('A9').Copy()
$ws.Range('A16').PasteSpecial(-4122)
$ws.Range('A16').Value = 'Syllabus:'
$ws.Range('B9').Copy()
$ws.Range('B16').PasteSpecial(-4122)
$ws.Range('B16').Value = '1.Fermentation submerged in bioreactor, involving setting, sterilization, medium addition, cultivation monitoring and samples analysis.
2.Solid state fermentation followed by recuperation and characterization of the produced enzymes.
3.Project of enzymes purification by software: the groups receive a problem- mixture and present logical sequences of purification with the respective results and their discussion.
4.Enzymes characterization in relation to molar mass: calibration of a chromatographic column with known proteins and determination of molar mass of problem-enzyme; determination of problem-enzyme by 280nm absorption and by specific activity.
5.Use of vegetal and microbial origin amylases for starch processing and ethanol fermentation.'
$ws.Range('C9').Copy()
$ws.Range('C16').PasteSpecial(-4122)
$ws.Range('C16').Value = '1.Fermentation submerged in bioreactor, involving setting, sterilization, medium addition, cultivation monitoring and samples analysis.
2.Solid state fermentation followed by recuperation and characterization of the produced enzymes.
3.Project of enzymes purification by software: the groups receive a problem- mixture and present logical sequences of purification with the respective results and their discussion.
4.Enzymes characterization in relation to molar mass: calibration of a chromatographic column with known proteins and determination of molar mass of problem-enzyme; determination of problem-enzyme by 280nm absorption and by specific activity.
5.Use of vegetal and microbial origin amylases for starch processing and ethanol fermentation.'
$ws.Rows.Item(16).RowHeight = 120

# Row 17
$ws.Range('A9').Copy()
$ws.Range('A17').PasteSpecial(-4122)
$ws.Range('A17').Value = 'Avaliação:'

# Row 18
$ws.Range('A9').Copy()
$ws.Range('A18').PasteSpecial(-4122)
$ws.Range('A18').Value = 'Método:'
$ws.Range('B9').Copy()
$ws.Range('B18').PasteSpecial(-4122)
$ws.Range('B18').Value = '3444370 - Rita de Cássia Lacerda Brambilla Rodrigues'
$ws.Range('C9').Copy()
$ws.Range('C18').PasteSpecial(-4122)
$ws.Range('C18').Value = '3444370 - Rita de Cássia Lacerda Brambilla Rodrigues'
$ws.Rows.Item(18).RowHeight = 60

# Row 19
$ws.Range('A9').Copy()
$ws.Range('A19').PasteSpecial(-4122)
$ws.Range('A19').Value = 'Critério:'
$ws.Range('B9').Copy()
$ws.Range('B19').PasteSpecial(-4122)
$ws.Range('B19').Value = 'Relatórios e seminários sobre os experimentos'
$ws.Range('C9').Copy()
$ws.Range('C19').PasteSpecial(-4122)
$ws.Range('C19').Value = 'Relatórios e seminários sobre os experimentos'
$ws.Rows.Item(19).RowHeight = 60

# Row 20
$ws.Range('A9').Copy()
$ws.Range('A20').PasteSpecial(-4122)
$ws.Range('A20').Value = 'Norma de recuperação:'
$ws.Range('B9').Copy()
$ws.Range('B20').PasteSpecial(-4122)
$ws.Range('B20').Value = 'Média aritmética entre os relatórios e seminários'
$ws.Range('C9').Copy()
$ws.Range('C20').PasteSpecial(-4122)
$ws.Range('C20').Value = 'Média aritmética entre os relatórios e seminários'
$ws.Rows.Item(20).RowHeight = 60

# Row 21
$ws.Range('A9').Copy()
$ws.Range('A21').PasteSpecial(-4122)
$ws.Range('A21').Value = 'Bibliografia:'
$ws.Range('B9').Copy()
$ws.Range('B21').PasteSpecial(-4122)
$ws.Range('B21').Value = 'A recuperação será feita por meio de prova escrita (PR) e a média final (MF) será calculada pela equação: MF = (NF + PR)/2.'
$ws.Range('C9').Copy()
$ws.Range('C21').PasteSpecial(-4122)
$ws.Range('C21').Value = 'A recuperação será feita por meio de prova escrita (PR) e a média final (MF) será calculada pela equação: MF = (NF + PR)/2.'
$ws.Rows.Item(21).RowHeight = 120

# Row 22
$ws.Range('A9').Copy()
$ws.Range('A22').PasteSpecial(-4122)
$ws.Range('A22').Value = 'Requisitos:'

# Row 23
$ws.Range('B9').Copy()
$ws.Range('B23').PasteSpecial(-4122)
$ws.Range('B23').Value = 'LOT2013 -  Engenharia Bioquímica I  (Requisito fraco)
'
$ws.Range('C9').Copy()
$ws.Range('C23').PasteSpecial(-4122)
$ws.Range('C23').Value = 'LOT2013 -  Engenharia Bioquímica I  (Requisito fraco)
'
$ws.Rows.Item(23).RowHeight = 30

# Row 24
$ws.Range('B9').Copy()
$ws.Range('B24').PasteSpecial(-4122)
$ws.Range('B24').Value = 'LOT2017 -  Enzimologia  (Requisito fraco)
'
$ws.Range('C9').Copy()
$ws.Range('C24').PasteSpecial(-4122)
$ws.Range('C24').Value = 'LOT2017 -  Enzimologia  (Requisito fraco)
'
$ws.Rows.Item(24).RowHeight = 30

$excel.CutCopyMode = $false